$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the historical_growth_revenue_last_5_years value in D2/D3
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 7).Value = -24.02547770700637   # G
    $ws.Cells.Item($r, 8).Value = -27.07006369426752   # H
    $ws.Cells.Item($r, 9).Value = -42.1656050955414    # I
    $ws.Cells.Item($r, 10).Value = -42.1656050955414   # J
    $ws.Cells.Item($r, 11).Value = -8.4                # K
    $ws.Cells.Item($r, 12).Value = -53.50318471337579  # L

    $ws.Cells.Item($r, 21).Value = 5.3                 # U
    $ws.Cells.Item($r, 22).Value = 0.1009523809523809  # V
    $ws.Cells.Item($r, 23).Value = -1.031941031941032  # W
    $ws.Cells.Item($r, 24).Value = 0.0589959088253032  # X
    $ws.Cells.Item($r, 25).Value = -1.090936940766335  # Y
    $ws.Cells.Item($r, 26).Value = 0.0215363511659808  # Z
    $ws.Cells.Item($r, 27).Value = -0.9080932784636488 # AA
    $ws.Cells.Item($r, 28).Value = 0.05193998712349561 # AB
    $ws.Cells.Item($r, 29).Value = -0.9600332655871444 # AC
    $ws.Cells.Item($r, 30).Value = 11.7                # AD
    $ws.Cells.Item($r, 32).Value = 11.7                # AF
    $ws.Cells.Item($r, 33).Value = 6.399999999999999   # AG
    $ws.Cells.Item($r, 34).Value = 0.1822429906542056  # AH
    $ws.Cells.Item($r, 35).Value = 0.9503695881731784  # AI
    $ws.Cells.Item($r, 36).Value = 0.1086587436332767  # AJ
    $ws.Cells.Item($r, 37).Value = 0.9128512337754957  # AK
    $ws.Cells.Item($r, 38).Value = 0.731               # AL
    $ws.Cells.Item($r, 39).Value = 0.728               # AM
    $ws.Cells.Item($r, 40).Value = -1.783536585365854  # AN
    $ws.Cells.Item($r, 41).Value = -9.056087551299591  # AO
    $ws.Cells.Item($r, 42).Value = -0.975609756097561  # AP
    $ws.Cells.Item($r, 43).Value = -9.093406593406593  # AQ
}
